$wb = $excel.ActiveWorkbook

# Update "想去人数" (number of people interested) for two events.
# These values are duplicated across the "展览" and "全部类型" sheets.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F4").Value = 85
    $ws.Range("F5").Value = 22
}
